# Auto-generated edit script: updates crypto price/volume data per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on columns B:E so numeric-looking strings (e.g. "61.815.26")
# are not reinterpreted by Excel as numbers/dates.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '61.815.26'
$ws.Range("E2").Value = '  +4.21%  '

$ws.Range("D3").Value = '3.071.81'
$ws.Range("E3").Value = '  +2.65%  '

$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").Value = '579.30'
$ws.Range("E5").Value = '  +2.98%  '

$ws.Range("D6").Value = '141.64'
$ws.Range("E6").Value = '  +2.93%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '3.061.41'
$ws.Range("E8").Value = '  +2.78%  '

$ws.Range("E9").Value = '  +1.05%  '

$ws.Range("E10").Value = '  +5.47%  '

$ws.Range("D11").Value = '5.75'
$ws.Range("E11").Value = '  +11.75%  '

$ws.Range("D12").Value = '0.465'
$ws.Range("E12").Value = '  +2.02%  '

$ws.Range("E13").Value = '  +4.35%  '

$ws.Range("D14").Value = '35.19'
$ws.Range("E14").Value = '  +4.54%  '

$ws.Range("E15").Value = '  +0.21%  '

$ws.Range("D16").Value = '3.582.42'
$ws.Range("E16").Value = '  +2.50%  '

$ws.Range("D17").Value = '7.26'
$ws.Range("E17").Value = '  +0.51%  '

$ws.Range("D18").Value = '3.071.86'
$ws.Range("E18").Value = '  +2.55%  '

$ws.Range("D19").Value = '61.761.61'
$ws.Range("E19").Value = '  +4.15%  '

$ws.Range("D20").Value = '446.56'
$ws.Range("E20").Value = '  +3.99%  '

$ws.Range("E21").Value = '  +2.16%  '

$ws.Range("D22").Value = '0.730'
$ws.Range("E22").Value = '  +1.69%  '

$ws.Range("D23").Value = '7.42'
$ws.Range("E23").Value = '  +4.29%  '

$ws.Range("D24").Value = '13.72'
$ws.Range("E24").Value = '  +3.01%  '

$ws.Range("D25").Value = '81.62'
$ws.Range("E25").Value = '  +0.92%  '

$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("E27").Value = '  +4.74%  '

$ws.Range("E28").Value = '  -0.10%  '

$ws.Range("D29").Value = '2.66'
$ws.Range("E29").Value = '  +4.75%  '

$ws.Range("D30").Value = '8.21'
$ws.Range("E30").Value = '  +6.26%  '

$ws.Range("E31").Value = '  +11.02%  '

$ws.Range("E32").Value = '  +12.91%  '

$ws.Range("D33").Value = '26.79'
$ws.Range("E33").Value = '  +4.17%  '

$ws.Range("E34").Value = '  +4.90%  '

$ws.Range("D35").Value = '0.0₃0790'
$ws.Range("E35").Value = '  +3.36%  '

$ws.Range("D36").Value = '6.04'
$ws.Range("E36").Value = '  +2.30%  '

$ws.Range("D37").Value = '2.18'
$ws.Range("E37").Value = '  +5.04%  '

$ws.Range("E38").Value = '  +2.07%  '

$ws.Range("D39").Value = '2.96'
$ws.Range("E39").Value = '  +8.67%  '

$ws.Range("E40").Value = '  +0.99%  '

$ws.Range("D41").Value = '420.34'
$ws.Range("E41").Value = '  +4.77%  '

$ws.Range("D42").Value = '2.959.90'
$ws.Range("E42").Value = '  +7.34%  '

$ws.Range("D43").Value = '0.0370'
$ws.Range("E43").Value = '  +5.17%  '

$ws.Range("D44").Value = '0.276'
$ws.Range("E44").Value = '  +9.89%  '

$ws.Range("E45").Value = '  +0.55%  '

$ws.Range("D46").Value = '2.12'
$ws.Range("E46").Value = '  +6.28%  '

$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").Value = '35.10'
$ws.Range("E48").Value = '  +0.95%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = '124.24'
$ws.Range("E49").Value = '  +3.34%  '

$ws.Range("E50").Value = '  +0.19%  '

$ws.Range("D51").Value = '24.29'
$ws.Range("E51").Value = '  +3.79%  '

